# Add 10 new commodity/group rows to the eNamMandiDataMaster table
# (VegMaster sheet, Table2), growing it from A1:B423 to A1:B432.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

$newRows = @(
    @("COPRA GRADE-1", "Misc"),
    @("SAL SEED", "Oil Seeds"),
    @("GROUNDNUT TJ", "Oil Seeds"),
    @("APPLE -RICH RED SMALL", "Fruits"),
    @("NAVANE", "Food grains/Cereals"),
    @("APPLE-RED ROYAL SMALL", "Fruits"),
    @("APPLE-GREEN SMITH XL", "Fruits"),
    @("PADDY CO 72", "Food grains/Cereals"),
    @("PADDY JAWAFOOL", "Food grains/Cereals"),
    @("WET-BULB TURMERIC", "Spices")
)

$startRow = $ws.UsedRange.Rows.Count + 1
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value2 = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value2 = $newRows[$i][1]
}

$lastRow = $startRow + $newRows.Count - 1
$lastCol = $lo.Range.Columns.Count
$lo.Resize($ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol)))

# Match the recorded view state from the edit: scrolled/selected near the
# newly appended rows.
$ws.Range("A395").Select()
$ws.Range("B428").Select()
